# Commit: "Change names from *img to img*"
#
# Renames the seven "*img" sheet tabs to "img*" (himg->imgh, timg->imgt,
# simg->imgs, gimg->imgg, wimg->imgw, bimg->imgb, eimg->imge), leaving
# every other sheet name untouched.
#
# The rename also moves which tab is shown as active/selected: the
# "himg" (-> "imgh") tab ends up the selected one, and "holiday" (which
# used to be the selected tab) is no longer selected.

$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $wb.Worksheets.Item($oldName).Name = $renames[$oldName]
}

# Make "imgh" (previously "himg") the active/selected sheet, which also
# clears the previous selection on "holiday".
$wb.Worksheets.Item("imgh").Activate()
